# Updates cryptos list values (price & 1h volume change) and fixes the
# VeChain/Maker row ordering, per the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row=2; D="68.253.18"; E="  +0.05%  " },
    @{ Row=3; D="3.721.55"; E="  -2.57%  " },
    @{ Row=4; D="1.00"; E="  +0.18%  " },
    @{ Row=5; D="598.67"; E="  +0.46%  " },
    @{ Row=6; D="167.84"; E="  -2.93%  " },
    @{ Row=7; D="3.720.48"; E="  -2.61%  " },
    @{ Row=8; E="  +0.03%  " },
    @{ Row=9; E="  +0.04%  " },
    @{ Row=10; D="0.169"; E="  +5.68%  " },
    @{ Row=11; D="6.24"; E="  -1.51%  " },
    @{ Row=12; D="0.462"; E="  -1.99%  " },
    @{ Row=13; D="38.36"; E="  -0.31%  " },
    @{ Row=14; D="0.0000247"; E="  +0.83%  " },
    @{ Row=15; D="4.342.70"; E="  -2.36%  " },
    @{ Row=16; D="3.726.16"; E="  -2.20%  " },
    @{ Row=17; D="68.235.96"; E="  -0.11%  " },
    @{ Row=18; D="7.33"; E="  +1.01%  " },
    @{ Row=20; D="17.23"; E="  +7.36%  " },
    @{ Row=21; D="491.07"; E="  +0.18%  " },
    @{ Row=22; D="9.27"; E="  -1.07%  " },
    @{ Row=23; D="0.725"; E="  -1.95%  " },
    @{ Row=24; D="84.98"; E="  -1.30%  " },
    @{ Row=25; D="0.0000143"; E="  +2.62%  " },
    @{ Row=26; D="2.32"; E="  -2.18%  " },
    @{ Row=27; D="12.33"; E="  +0.38%  " },
    @{ Row=28; D="10.12"; E="  -0.65%  " },
    @{ Row=29; E="  +0.05%  " },
    @{ Row=30; D="2.92"; E="  -0.62%  " },
    @{ Row=31; D="7.90"; E="  +3.45%  " },
    @{ Row=32; D="2.38"; E="  -2.46%  " },
    @{ Row=33; D="31.59"; E="  -4.59%  " },
    @{ Row=34; D="3.863.24"; E="  -2.24%  " },
    @{ Row=35; E="  -2.02%  " },
    @{ Row=36; D="3.667.06"; E="  -2.40%  " },
    @{ Row=37; D="1.00"; E="  +0.33%  " },
    @{ Row=38; E="  -0.50%  " },
    @{ Row=39; D="5.85"; E="  +0.25%  " },
    @{ Row=40; E="  -3.18%  " },
    @{ Row=41; D="0.324"; E="  -0.56%  " },
    @{ Row=42; D="431.65"; E="  -5.60%  " },
    @{ Row=43; D="48.77"; E="  -0.97%  " },
    @{ Row=44; E="  -2.11%  " },
    @{ Row=45; D="2.86"; E="  -1.54%  " },
    @{ Row=46; D="8.44"; E="  +1.11%  " },
    @{ Row=48; D="40.48"; E="  -2.40%  " },
    @{ Row=49; D="141.39"; E="  +2.66%  " },
    @{ Row=50; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="2.761.32"; E="  -3.20%  " },
    @{ Row=51; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0352"; E="  -0.31%  " }
)

foreach ($u in $rowUpdates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Force text storage so numeric-looking price strings (e.g. "1.00")
        # keep their original formatting instead of being coerced to numbers.
        $dCell = $ws.Range("D" + $u.Row)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $u.Row).Value = $u.E }
}
